$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6, column B: "B0,I, biomass in initial year" -> "N0,I, abundance in initial year" ---
$cell = $ws.Range("B6")
$cell.Value = "N0,I, abundance in initial year"
$sub = $cell.Characters(2, 3)
$sub.Font.Subscript = $true
$rest = $cell.Characters(5, 27)
$rest.Font.Name = "Calibri"
$rest.Font.Size = 11
$rest.Font.Color = 0

# --- Row 7, column C: "Observed (observed catch / observed biomass)" -> "...observed abundance)" ---
$ws.Range("C7").Value = "Observed (observed catch / observed abundance)"

# --- Row 10, column B: "Bi,t, biomass time series" -> "Ni,t, abundance time series" ---
$cell = $ws.Range("B10")
$cell.Value = "Ni,t, abundance time series"
$sub = $cell.Characters(2, 3)
$sub.Font.Subscript = $true
$rest = $cell.Characters(5, 23)
$rest.Font.Name = "Calibri"
$rest.Font.Size = 11
$rest.Font.Color = 0

# --- Column C width: widen slightly to fit the updated text ---
$ws.Columns("C").ColumnWidth = 38

# --- Move the saved selection to C11 (matches the view state recorded in the diff) ---
$ws.Range("C11").Select()
